$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report year labels to include "CY " (calendar year)
$ws.Range("A2").Value = "report for CY 2021"
$ws.Range("A3").Value = "report for CY 2022"
$ws.Range("A4").Value = "report for CY 2023"

# Column A's best-fit width grows to accommodate the longer labels
$ws.Columns.Item(1).ColumnWidth = 15.29

# Move the active selection to B10 (matches the saved selection state)
$ws.Range("B10").Select()
